$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.259004950523376
$ws.Range("B1").Value = 2.954211473464966
$ws.Range("C1").Value = 6.191031455993652
$ws.Range("D1").Value = 4.909485816955566
$ws.Range("E1").Value = 1.240481734275818
